$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 3.071615738262478
$ws.Range("B3").Value = 10.00000009988942
$ws.Range("B5").Value = 1.37352946067567
$ws.Range("B7").Value = 0.09711556084068254
$ws.Range("B8").Value = 3.606375541929839
$ws.Range("B9").Value = 1.039755944406956
$ws.Range("B10").Value = 0.8327578712298451
$ws.Range("B11").Value = -0.2069980731771111
$ws.Range("B12").Value = 0.1000000097260039
$ws.Range("B13").Value = 0.7000000099810222
$ws.Range("B14").Value = 0.03545639447291898
$ws.Range("B15").Value = 57393.16471833006
$ws.Range("B16").Value = 0.633625850059178
$ws.Range("B17").Value = 0.172939035495789
$ws.Range("B18").Value = 0.1652042599136138
$ws.Range("B19").Value = 0.01942631860822956
$ws.Range("B20").Value = 0.0000000000000256873097489645
$ws.Range("B21").Value = 0.00000000000007994990901557821
$ws.Range("B22").Value = -3.063945696245146
$ws.Range("B23").Value = 1.377461840494002
$ws.Range("B24").Value = -20.24761807578205
$ws.Range("B25").Value = 0.4914338153123421
$ws.Range("B26").Value = 3.064047296725472
$ws.Range("B27").Value = 0.04284846721856334
$ws.Range("B28").Value = 1.643643111380612
$ws.Range("B29").Value = 1.377461840494002
$ws.Range("B30").Value = 0.00009387763229539152
$ws.Range("B31").Value = 0.0000000008060675216719966
$ws.Range("B32").Value = 0.02637643377149679
$ws.Range("B33").Value = 0
$ws.Range("B34").Value = 0.007519776023409953
$ws.Range("B35").Value = 0.002229022651039306
$ws.Range("B36").Value = 0.0008196617250040024
$ws.Range("B37").Value = 0.002990680391247099
$ws.Range("B38").Value = 0
$ws.Range("B39").Value = -0.00004745381095122977
$ws.Range("B40").Value = 0
$ws.Range("B41").Value = 97.11556084068255
$ws.Range("B42").Value = 0.2141031623628998
$ws.Range("B43").Value = 0.03424192902816622
$ws.Range("B44").Value = 0.01635522173144777
$ws.Range("B45").Value = 0.001923205542214726
$ws.Range("B50").Value = 0.002601469290930152
$ws.Range("B52").Value = -0.000000009999534314054936
$ws.Range("B53").Value = 0.007193745247457965
$ws.Range("B54").Value = -0.1156359638013393
$ws.Range("B55").Value = -0.3224945337161784
$ws.Range("B56").Value = -0.0000000000000001013084226683912
$ws.Range("B57").Value = -0.0000000000000001073207592805038
$ws.Range("B58").Value = 1.163036398800098
$ws.Range("B59").Value = 0.0000000000000004364676367446859
$ws.Range("B60").Value = 0.3224945337161784
$ws.Range("B61").Value = -0.0000000000000001013084226683912
$ws.Range("B62").Value = 0.000000000000001688388232663155
$ws.Range("B63").Value = -0.00000000000000004348700486638601
$ws.Range("B64").Value = 1.163752728361828
$ws.Range("B65").Value = 0.0000000000000004367364629579164
$ws.Range("B66").Value = 0.3226931623818295
$ws.Range("B67").Value = -0.0000000000000001599852970987496
$ws.Range("B68").Value = 0.000000000000006189842183794375
$ws.Range("B69").Value = -0.0000000000000000686742642934687
$ws.Range("B70").Value = 0.1834148783059018
$ws.Range("B71").Value = 0.1390796554102766
$ws.Range("B72").Value = 3.315227952113323
$ws.Range("B73").Value = 1.917688490849432
$ws.Range("B74").Value = 0.00000000000174043713339603
$ws.Range("B75").Value = 0.000000000006749252350959501
$ws.Range("B76").Value = -0.3315228274552523
$ws.Range("B77").Value = 0.000000000002212146399692676
$ws.Range("B78").Value = 0.05366999389056305
$ws.Range("B79").Value = -0.00002281246484170344
$ws.Range("B80").Value = -0.000004248266888030819
$ws.Range("B81").Value = -0.07732045567273302
$ws.Range("B82").Value = -0.07999999002558271
$ws.Range("B83").Value = -0.000002743322328451323
$ws.Range("B84").Value = 0.02999999000734206
$ws.Range("B85").Value = 0.2722204216888156
$ws.Range("B86").Value = -0.04135978094277348
$ws.Range("B87").Value = -0.006731097246748252
$ws.Range("B88").Value = -0.2086574852286462
$ws.Range("B89").Value = -0.7434467110570215
$ws.Range("B90").Value = -0.005971531420664135
$ws.Range("B91").Value = -0.1173659766867269
$ws.Range("B92").Value = -0.7768111938772915
$ws.Range("B93").Value = -0.2135889760795928
$ws.Range("B94").Value = 0.00000000000001422554333503584
$ws.Range("B95").Value = 0.00000000000006015570153565853
$ws.Range("B96").Value = 0.01683213031752181
$ws.Range("B97").Value = 0.00000000000003697139921994185
$ws.Range("B98").Value = -0.001289541961213558
$ws.Range("B99").Value = -0.000114908025150573
$ws.Range("B100").Value = 0.10172508052207
$ws.Range("B101").Value = 0.4110248747055887
$ws.Range("B102").Value = -0.0002170000508277858
$ws.Range("B103").Value = -0.02454754798352729
